$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numericTextCells = @("D5","D6","D7","D11","D13","D14","D18","D20","D21","D22","D23","D24","D25","D28","D29","D34","D35","D37","D38","D39","D40","D41","D42","D44","D46","D49")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.245.05"
$ws.Range("E2").Value = "  -4.75%  "
$ws.Range("D3").Value = "3.259.03"
$ws.Range("E3").Value = "  -7.22%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "596.59"
$ws.Range("E5").Value = "  -4.30%  "
$ws.Range("D6").Value = "150.71"
$ws.Range("E6").Value = "  -12.64%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "3.249.48"
$ws.Range("E8").Value = "  -7.35%  "
$ws.Range("E9").Value = "  -11.40%  "
$ws.Range("E10").Value = "  -13.96%  "
$ws.Range("D11").Value = "6.66"
$ws.Range("E11").Value = "  -6.79%  "
$ws.Range("E12").Value = "  -14.00%  "
$ws.Range("D13").Value = "38.22"
$ws.Range("E13").Value = "  -17.55%  "
$ws.Range("D14").Value = "0.0000243"
$ws.Range("E14").Value = "  -12.04%  "
$ws.Range("D15").Value = "3.780.49"
$ws.Range("E15").Value = "  -7.51%  "
$ws.Range("D16").Value = "67.253.95"
$ws.Range("E16").Value = "  -4.88%  "
$ws.Range("D17").Value = "3.259.06"
$ws.Range("E17").Value = "  -7.46%  "
$ws.Range("D18").Value = "535.28"
$ws.Range("E18").Value = "  -12.15%  "
$ws.Range("E19").Value = "  -6.18%  "
$ws.Range("D20").Value = "7.21"
$ws.Range("E20").Value = "  -14.05%  "
$ws.Range("D21").Value = "15.09"
$ws.Range("E21").Value = "  -14.94%  "
$ws.Range("D22").Value = "0.761"
$ws.Range("E22").Value = "  -13.66%  "
$ws.Range("D23").Value = "7.87"
$ws.Range("E23").Value = "  -13.40%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "13.57"
$ws.Range("E24").Value = "  -12.89%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "85.15"
$ws.Range("E25").Value = "  -12.49%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  -12.69%  "
$ws.Range("D28").Value = "29.30"
$ws.Range("E28").Value = "  -12.55%  "
$ws.Range("D29").Value = "8.02"
$ws.Range("E29").Value = "  -11.42%  "
$ws.Range("E30").Value = "  -16.93%  "
$ws.Range("E31").Value = "  -10.70%  "
$ws.Range("E32").Value = "  -11.35%  "
$ws.Range("E33").Value = "  -17.72%  "
$ws.Range("D34").Value = "541.46"
$ws.Range("E34").Value = "  -13.65%  "
$ws.Range("D35").Value = "5.71"
$ws.Range("E35").Value = "  -16.25%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").Value = "0.0452"
$ws.Range("E37").Value = "  -8.29%  "
$ws.Range("D38").Value = "53.33"
$ws.Range("E38").Value = "  -5.85%  "
$ws.Range("D39").Value = "0.0855"
$ws.Range("E39").Value = "  -13.86%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "0.128"
$ws.Range("E40").Value = "  -9.78%  "
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").Value = "9.09"
$ws.Range("E41").Value = "  -15.86%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "2.73"
$ws.Range("E42").Value = "  -20.38%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.932.76"
$ws.Range("E43").Value = "  -12.29%  "
$ws.Range("D44").Value = "0.263"
$ws.Range("E44").Value = "  -15.72%  "
$ws.Range("D45").Value = "0.0₃0583"
$ws.Range("E45").Value = "  -19.11%  "
$ws.Range("D46").Value = "26.62"
$ws.Range("E46").Value = "  -16.78%  "
$ws.Range("E47").Value = "  -14.27%  "
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").Value = "127.65"
$ws.Range("E49").Value = "  -3.96%  "
$ws.Range("E50").Value = "  -21.97%  "
$ws.Range("E51").Value = "  -12.56%  "
